# Exponential Growth.xlsx - add "Phase one monochrome conversion" (Half Turns)
# columns to the Sheet1 summary table, plus new "Total" columns for both the
# Quarter Turns and Half Turns blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Break apart the merged header cells that need to move / resize before we
#    touch any values, then re-merge them in their new layout.
# ---------------------------------------------------------------------------
$ws.Range("H4:I4").UnMerge()
$ws.Range("J4:K4").UnMerge()
$ws.Range("G4:G5").UnMerge()
$ws.Range("L4:L5").UnMerge()
$ws.Range("H3:L3").UnMerge()
$ws.Range("C3:G3").UnMerge()

# ---------------------------------------------------------------------------
# 2. Row 3 - top banner ("Quarter Turns" now spans C:H, "Half Turns" moves to
#    I:N).
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = $null
$ws.Range("I3").Value = "Half Turns"
$ws.Range("J3").Value = $null
$ws.Range("K3").Value = $null
$ws.Range("L3").Value = $null
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = $null

# ---------------------------------------------------------------------------
# 3. Row 4 - phase banners, plus new "Total" banners in G and M.
#    G4:H4, I4:J4 and K4:L4 already have a pre-existing (non-default) style
#    on their top-left cell, so merging them later naturally keeps the header
#    look. M4 is a brand new cell (outside the old B3:L10 range) so it starts
#    out with the plain column style - give it the same center/center
#    alignment as the rest of the banner row *before* merging so the merged
#    M4:N4 range matches its neighbours instead of falling back to default.
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = "Total"
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = "Phase 1"
$ws.Range("J4").Value = $null
$ws.Range("K4").Value = "Phase 2"
$ws.Range("L4").Value = $null
$ws.Range("M4").Value = "Total"
$ws.Range("N4").Value = $null
$ws.Range("M4").HorizontalAlignment = -4108
$ws.Range("M4").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Row 5 - column captions (Depth / Count / Time). G5 and L5 are
#    pre-existing cells that used to be blank halves of a vertical merge
#    (G4:G5 / L4:L5) - now that the merge is gone they need the same
#    center-aligned caption style as their C5:F5/H5:K5 neighbours. M5/N5 are
#    brand new cells so need the same explicit alignment.
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = "Count"
$ws.Range("H5").Value = "Time"
$ws.Range("I5").Value = "Depth"
$ws.Range("J5").Value = "Count"
$ws.Range("K5").Value = "Depth"
$ws.Range("L5").Value = "Count"
$ws.Range("M5").Value = "Count"
$ws.Range("N5").Value = "Time"
$ws.Range("G5").HorizontalAlignment = -4108
$ws.Range("G5").VerticalAlignment = -4108
$ws.Range("L5").HorizontalAlignment = -4108
$ws.Range("L5").VerticalAlignment = -4108
$ws.Range("M5:N5").HorizontalAlignment = -4108
$ws.Range("M5:N5").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Data rows 6-10. Column G becomes a "Total" formula (Quarter Turns count),
#    the old Growth-Factor value moves from G to H ("Time"), and the whole
#    "Half Turns" block is populated in I:N (with M also a "Total" formula).
# ---------------------------------------------------------------------------

# Row 6
$ws.Range("H6").Value = 0.15
$ws.Range("G6").Formula = "=F6+D6"
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 480
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 23
$ws.Range("M6").Formula = "=L6+J6"
$ws.Range("N6").Value = 0.503

# Row 7
$ws.Range("H7").Value = 0.218
$ws.Range("G7").Formula = "=F7+D7"
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 180
$ws.Range("M7").Formula = "=L7+J7"
$ws.Range("N7").Value = 0.316

# Row 8
$ws.Range("H8").Value = 180.63300000000001
$ws.Range("G8").Formula = "=F8+D8"
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 246
$ws.Range("K8").Value = 6
$ws.Range("L8").Value = 451204
$ws.Range("M8").Formula = "=L8+J8"
$ws.Range("N8").Value = 429.428

# Row 9
$ws.Range("H9").Value = 0.187
$ws.Range("G9").Formula = "=F9+D9"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 180
$ws.Range("M9").Formula = "=L9+J9"
$ws.Range("N9").Value = 0.22800000000000001

# Row 10
$ws.Range("H10").Value = 13324.681
$ws.Range("G10").Formula = "=F10+D10"
$ws.Range("M10").Formula = "=L10+J10"

# ---------------------------------------------------------------------------
# 6. Re-merge the header cells into their new layout.
# ---------------------------------------------------------------------------
$ws.Range("C3:H3").Merge()
$ws.Range("I3:N3").Merge()
$ws.Range("G4:H4").Merge()
$ws.Range("I4:J4").Merge()
$ws.Range("K4:L4").Merge()
$ws.Range("M4:N4").Merge()

# ---------------------------------------------------------------------------
# 7. Selection moves to C10:D10.
# ---------------------------------------------------------------------------
$ws.Range("C10:D10").Select()
